$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing board-token values (BNA -> BNAKU, RPA -> RPALT)
$ws.Range("A2").Value = "BNAKU"
$ws.Range("A3").Value = "RPALT"

# Move selection to A2
$ws.Range("A2").Select()
